{"js": "// Replace the date line and all 100 arithmetic-problem answers in the table\n// with the new values from the diff. Every \"old\" text below occurs exactly\n// once in the document (verified against before.docx), so a targeted\n// body.search(...) + insertText(..., replace) per pair is safe and avoids\n// depending on paragraph/table traversal order.\nconst replacements = [\n  [\"2024-07-27 Saturday\", \"2024-07-28 Sunday\"],\n  [\"10+63=73\", \"94-88=6\"],\n  [\"51-45=6\", \"2+75=77\"],\n  [\"88-4=84\", \"60-9=51\"],\n  [\"46+36=82\", \"15+53=68\"],\n  [\"31-19=12\", \"75+16=91\"],\n  [\"61-54=7\", \"29-18=11\"],\n  [\"37+8=45\", \"63+13=76\"],\n  [\"30+3=33\", \"76-71=5\"],\n  [\"27+40=67\", \"99-40=59\"],\n  [\"51+45=96\", \"99-49=50\"],\n  [\"6+13=19\", \"6+18=24\"],\n  [\"22-5=17\", \"68-32=36\"],\n  [\"3+32=35\", \"58-35=23\"],\n  [\"13+50=63\", \"7+3=10\"],\n  [\"22-1=21\", \"1+78=79\"],\n  [\"83-14=69\", \"14-10=4\"],\n  [\"71+16=87\", \"3+3=6\"],\n  [\"31+42=73\", \"22+77=99\"],\n  [\"13+38=51\", \"78-8=70\"],\n  [\"31-26=5\", \"0+87=87\"],\n  [\"66+28=94\", \"69-44=25\"],\n  [\"33+14=47\", \"47+42=89\"],\n  [\"27+2=29\", \"2+23=25\"],\n  [\"62+26=88\", \"43+34=77\"],\n  [\"66-27=39\", \"28-3=25\"],\n  [\"4+12=16\", \"85-22=63\"],\n  [\"12+17=29\", \"70+25=95\"],\n  [\"18+6=24\", \"4+70=74\"],\n  [\"27-15=12\", \"35-7=28\"],\n  [\"2+13=15\", \"77+14=91\"],\n  [\"75-48=27\", \"8-4=4\"],\n  [\"81-41=40\", \"81-81=0\"],\n  [\"65-52=13\", \"76-43=33\"],\n  [\"91-81=10\", \"3+21=24\"],\n  [\"87-45=42\", \"18+29=47\"],\n  [\"46+30=76\", \"28+23=51\"],\n  [\"62-12=50\", \"53+1=54\"],\n  [\"16-13=3\", \"89-8=81\"],\n  [\"74-36=38\", \"82+17=99\"],\n  [\"55-14=41\", \"26+28=54\"],\n  [\"11+84=95\", \"47+39=86\"],\n  [\"52-19=33\", \"44+44=88\"],\n  [\"43+1=44\", \"57-38=19\"],\n  [\"71-49=22\", \"86-26=60\"],\n  [\"96-33=63\", \"10+36=46\"],\n  [\"92-79=13\", \"37-9=28\"],\n  [\"40+31=71\", \"55+19=74\"],\n  [\"85-52=33\", \"43+52=95\"],\n  [\"10-3=7\", \"88-60=28\"],\n  [\"21+1=22\", \"65-4=61\"],\n  [\"56+9=65\", \"76+8=84\"],\n  [\"79-14=65\", \"89-81=8\"],\n  [\"65+22=87\", \"67-33=34\"],\n  [\"26-0=26\", \"95-84=11\"],\n  [\"93-38=55\", \"65-47=18\"],\n  [\"50+16=66\", \"29+31=60\"],\n  [\"12+63=75\", \"21-20=1\"],\n  [\"13+45=58\", \"25+7=32\"],\n  [\"5+86=91\", \"69-51=18\"],\n  [\"35+12=47\", \"86-76=10\"],\n  [\"54-11=43\", \"56-7=49\"],\n  [\"27+48=75\", \"96-26=70\"],\n  [\"86-59=27\", \"88-53=35\"],\n  [\"17+21=38\", \"23-17=6\"],\n  [\"31+2=33\", \"17+18=35\"],\n  [\"43+53=96\", \"54+14=68\"],\n  [\"64-54=10\", \"44+3=47\"],\n  [\"84-37=47\", \"88-11=77\"],\n  [\"10+85=95\", \"57+16=73\"],\n  [\"43+51=94\", \"31+43=74\"],\n  [\"27-23=4\", \"33+56=89\"],\n  [\"57+11=68\", \"39-4=35\"],\n  [\"68-27=41\", \"76-46=30\"],\n  [\"75+14=89\", \"64+12=76\"],\n  [\"77-71=6\", \"89-57=32\"],\n  [\"98-5=93\", \"65-10=55\"],\n  [\"99-88=11\", \"94-56=38\"],\n  [\"69+6=75\", \"36+40=76\"],\n  [\"13+81=94\", \"18+29=47\"],\n  [\"74-64=10\", \"22+58=80\"],\n  [\"30+1=31\", \"85-54=31\"],\n  [\"80-3=77\", \"46-7=39\"],\n  [\"24+48=72\", \"92-58=34\"],\n  [\"44-23=21\", \"85+12=97\"],\n  [\"9+8=17\", \"71+2=73\"],\n  [\"28-25=3\", \"80-79=1\"],\n  [\"48+7=55\", \"19+28=47\"],\n  [\"51-22=29\", \"53-2=51\"],\n  [\"74-41=33\", \"82-79=3\"],\n  [\"52+24=76\", \"1+3=4\"],\n  [\"22+76=98\", \"86-76=10\"],\n  [\"29+32=61\", \"60-27=33\"],\n  [\"25-1=24\", \"75-57=18\"],\n  [\"15+37=52\", \"89-47=42\"],\n  [\"24+22=46\", \"41-27=14\"],\n  [\"32+22=54\", \"78-51=27\"],\n  [\"33+19=52\", \"44+29=73\"],\n  [\"22+12=34\", \"76-62=14\"],\n  [\"83-17=66\", \"17+25=42\"],\n  [\"13+46=59\", \"69-0=69\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found, cannot apply edit: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all 100 arithmetic-problem answers in the table\n# with the new values from the diff. Every \"old\" text below occurs exactly\n# once in the document (verified against before.docx), so Find/Replace per\n# pair, scoped to a fresh Range each time, is safe and avoids depending on\n# paragraph/table traversal order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('2024-07-27 Saturday', '2024-07-28 Sunday')\n    ,@('10+63=73', '94-88=6')\n    ,@('51-45=6', '2+75=77')\n    ,@('88-4=84', '60-9=51')\n    ,@('46+36=82', '15+53=68')\n    ,@('31-19=12', '75+16=91')\n    ,@('61-54=7', '29-18=11')\n    ,@('37+8=45', '63+13=76')\n    ,@('30+3=33', '76-71=5')\n    ,@('27+40=67', '99-40=59')\n    ,@('51+45=96', '99-49=50')\n    ,@('6+13=19', '6+18=24')\n    ,@('22-5=17', '68-32=36')\n    ,@('3+32=35', '58-35=23')\n    ,@('13+50=63', '7+3=10')\n    ,@('22-1=21', '1+78=79')\n    ,@('83-14=69', '14-10=4')\n    ,@('71+16=87', '3+3=6')\n    ,@('31+42=73', '22+77=99')\n    ,@('13+38=51', '78-8=70')\n    ,@('31-26=5', '0+87=87')\n    ,@('66+28=94', '69-44=25')\n    ,@('33+14=47', '47+42=89')\n    ,@('27+2=29', '2+23=25')\n    ,@('62+26=88', '43+34=77')\n    ,@('66-27=39', '28-3=25')\n    ,@('4+12=16', '85-22=63')\n    ,@('12+17=29', '70+25=95')\n    ,@('18+6=24', '4+70=74')\n    ,@('27-15=12', '35-7=28')\n    ,@('2+13=15', '77+14=91')\n    ,@('75-48=27', '8-4=4')\n    ,@('81-41=40', '81-81=0')\n    ,@('65-52=13', '76-43=33')\n    ,@('91-81=10', '3+21=24')\n    ,@('87-45=42', '18+29=47')\n    ,@('46+30=76', '28+23=51')\n    ,@('62-12=50', '53+1=54')\n    ,@('16-13=3', '89-8=81')\n    ,@('74-36=38', '82+17=99')\n    ,@('55-14=41', '26+28=54')\n    ,@('11+84=95', '47+39=86')\n    ,@('52-19=33', '44+44=88')\n    ,@('43+1=44', '57-38=19')\n    ,@('71-49=22', '86-26=60')\n    ,@('96-33=63', '10+36=46')\n    ,@('92-79=13', '37-9=28')\n    ,@('40+31=71', '55+19=74')\n    ,@('85-52=33', '43+52=95')\n    ,@('10-3=7', '88-60=28')\n    ,@('21+1=22', '65-4=61')\n    ,@('56+9=65', '76+8=84')\n    ,@('79-14=65', '89-81=8')\n    ,@('65+22=87', '67-33=34')\n    ,@('26-0=26', '95-84=11')\n    ,@('93-38=55', '65-47=18')\n    ,@('50+16=66', '29+31=60')\n    ,@('12+63=75', '21-20=1')\n    ,@('13+45=58', '25+7=32')\n    ,@('5+86=91', '69-51=18')\n    ,@('35+12=47', '86-76=10')\n    ,@('54-11=43', '56-7=49')\n    ,@('27+48=75', '96-26=70')\n    ,@('86-59=27', '88-53=35')\n    ,@('17+21=38', '23-17=6')\n    ,@('31+2=33', '17+18=35')\n    ,@('43+53=96', '54+14=68')\n    ,@('64-54=10', '44+3=47')\n    ,@('84-37=47', '88-11=77')\n    ,@('10+85=95', '57+16=73')\n    ,@('43+51=94', '31+43=74')\n    ,@('27-23=4', '33+56=89')\n    ,@('57+11=68', '39-4=35')\n    ,@('68-27=41', '76-46=30')\n    ,@('75+14=89', '64+12=76')\n    ,@('77-71=6', '89-57=32')\n    ,@('98-5=93', '65-10=55')\n    ,@('99-88=11', '94-56=38')\n    ,@('69+6=75', '36+40=76')\n    ,@('13+81=94', '18+29=47')\n    ,@('74-64=10', '22+58=80')\n    ,@('30+1=31', '85-54=31')\n    ,@('80-3=77', '46-7=39')\n    ,@('24+48=72', '92-58=34')\n    ,@('44-23=21', '85+12=97')\n    ,@('9+8=17', '71+2=73')\n    ,@('28-25=3', '80-79=1')\n    ,@('48+7=55', '19+28=47')\n    ,@('51-22=29', '53-2=51')\n    ,@('74-41=33', '82-79=3')\n    ,@('52+24=76', '1+3=4')\n    ,@('22+76=98', '86-76=10')\n    ,@('29+32=61', '60-27=33')\n    ,@('25-1=24', '75-57=18')\n    ,@('15+37=52', '89-47=42')\n    ,@('24+22=46', '41-27=14')\n    ,@('32+22=54', '78-51=27')\n    ,@('33+19=52', '44+29=73')\n    ,@('22+12=34', '76-62=14')\n    ,@('83-17=66', '17+25=42')\n    ,@('13+46=59', '69-0=69')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $matched = $find.Execute(\n        $oldText,  # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $true,     # Format\n        $newText,  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n\n    if (-not $matched) {\n        throw \"Text not found, cannot apply edit: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
